$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1899
$ws.Range("F5").Value = 786
$ws.Range("F8").Value = 965
$ws.Range("F9").Value = 1647
$ws.Range("F10").Value = 1298
$ws.Range("F11").Value = 1584
$ws.Range("F12").Value = 74
$ws.Range("F13").Value = 1583
$ws.Range("F14").Value = 354
$ws.Range("F15").Value = 1719
$ws.Range("F16").Value = 818
$ws.Range("F17").Value = 1148
$ws.Range("F18").Value = 388
$ws.Range("F19").Value = 59
$ws.Range("F20").Value = 118
$ws.Range("F21").Value = 1946
$ws.Range("F22").Value = 266
$ws.Range("F23").Value = 831
$ws.Range("F24").Value = 1018
$ws.Range("F26").Value = 1299
$ws.Range("F27").Value = 1093
$ws.Range("F28").Value = 95
$ws.Range("F29").Value = 584
$ws.Range("F30").Value = 1219
$ws.Range("F31").Value = 916
$ws.Range("F32").Value = 1206
$ws.Range("F33").Value = 1155
$ws.Range("F34").Value = 298
$ws.Range("F35").Value = 90
$ws.Range("F36").Value = 903
$ws.Range("F37").Value = 5
$ws.Range("F38").Value = 1724
$ws.Range("F40").Value = 126
$ws.Range("F41").Value = 2095
$ws.Range("F42").Value = 104
$ws.Range("F44").Value = 195

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 1508
$ws.Range("F8").Value = 2620
$ws.Range("F12").Value = 270
$ws.Range("F21").Value = 329
$ws.Range("F22").Value = 92566
$ws.Range("F23").Value = 32
$ws.Range("F31").Value = 234
$ws.Range("F37").Value = 17
$ws.Range("G37").Value = 128
$ws.Range("F38").Value = 192
$ws.Range("F44").Value = 146
$ws.Range("F45").Value = 70

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 274
$ws.Range("F5").Value = 2944
$ws.Range("F6").Value = 4718
$ws.Range("F7").Value = 155
$ws.Range("F9").Value = 608
$ws.Range("F10").Value = 802
$ws.Range("F12").Value = 447
$ws.Range("F13").Value = 1182
$ws.Range("F14").Value = 328
$ws.Range("F15").Value = 815

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1899
$ws.Range("F3").Value = 274
$ws.Range("F5").Value = 802
$ws.Range("F7").Value = 447
$ws.Range("F8").Value = 447
$ws.Range("F9").Value = 1182
$ws.Range("F11").Value = 965
$ws.Range("F12").Value = 1647
$ws.Range("F13").Value = 1298
$ws.Range("F14").Value = 1584
$ws.Range("F15").Value = 74
$ws.Range("F16").Value = 1583
$ws.Range("F17").Value = 270
$ws.Range("F19").Value = 1719
$ws.Range("F20").Value = 1148
$ws.Range("F22").Value = 815
$ws.Range("F23").Value = 815
$ws.Range("F24").Value = 1946
$ws.Range("F25").Value = 266
$ws.Range("F26").Value = 831
$ws.Range("F27").Value = 1018
$ws.Range("F29").Value = 1299
$ws.Range("F30").Value = 329
$ws.Range("F31").Value = 1093
$ws.Range("F32").Value = 95
$ws.Range("F33").Value = 1219
$ws.Range("F34").Value = 916
$ws.Range("F35").Value = 1206
$ws.Range("F38").Value = 1155
$ws.Range("F39").Value = 298
$ws.Range("F40").Value = 903
$ws.Range("F42").Value = 5
$ws.Range("F43").Value = 1724
$ws.Range("F45").Value = 126
$ws.Range("F46").Value = 2095
$ws.Range("F47").Value = 104
$ws.Range("F49").Value = 196
$ws.Range("F53").Value = 70
